# Automatische test-sync: 2025-08-06 19:30:50
$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# Append new row 4 with the new test-mail entry
$logs.Range("A4").Value = "Zou jij dit even op kunnen pakken?"
$logs.Range("B4").Value = "mailmind.test@zohomail.eu"
$logs.Range("C4").Value = "Testmail #1: Zou jij dit even op kunnen pakken?"
$logs.Range("D4").Value = "Planning / Afspraak"
$logs.Range("E4").Value = "Bedankt, we hebben dit doorgestuurd naar planning@bedrijf.nl."
$logs.Range("F4").Value = "2025-08-06 19:29:59"
$logs.Range("G4").Value = "Ja"
$logs.Range("H4").Value = "Ja"
$logs.Range("I4").Value = "Nee"
$logs.Range("J4").Value = "Nee"

# Update the dashboard count for "Planning / Afspraak"
$dashboard.Range("B2").Value = 3

# Extend the conditional-formatting ranges so they keep covering the data
# (they previously ended at row 3, now they must end at row 4)
$colsToExtend = @("D", "G", "H", "I", "J")
foreach ($col in $colsToExtend) {
    $oldRange = $logs.Range("$col`2:$col`3")
    $newRange = $logs.Range("$col`2:$col`4")
    for ($i = 1; $i -le $oldRange.FormatConditions.Count; $i++) {
        $fc = $oldRange.FormatConditions.Item($i)
        $fc.ModifyAppliesToRange($newRange)
    }
}
